# "remote clone work in excel"
# Append two new rows (4 and 5) of single-letter string values below the
# existing data, mirroring the first three rows with new shared-string
# entries, and leave the selection on the last cell entered (C5).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "z"
$ws.Range("B4").Value = "y"
$ws.Range("C4").Value = "x"

$ws.Range("A5").Value = "v"
$ws.Range("B5").Value = "u"
$ws.Range("C5").Value = "t"

$ws.Range("C5").Select() | Out-Null
